# ---------------------------------------------------------------------------
# Edit: trims the "Limitations" discussion at the end of section 2.
#
#  1. The closing paragraph of the limitations narrative is reworded:
#       "Understanding these limitations is crucial for interpreting the
#        study's findings accurately and recognizing the potential for
#        variance in real-world applications. Further research and refined
#        methodologies could help address these limitations in future
#        analyses."
#     becomes:
#       "Further research, data collection, and methodological refinements
#        could help address these limitations in future analyses."
#
#  2. The nine bulleted "limitations" list items that followed it (the
#     numId=1005 list, from "Only includes parcels where residential use
#     is allowed" through "Assumes all data to be concurrent") are removed
#     entirely, since that information is now redundant with the prose
#     paragraph above it (and the earlier prose list covering the same
#     points).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Step 1: reword the closing "limitations" paragraph --------------------

$rightQuote = [char]0x2019
$oldParagraphText = "Understanding these limitations is crucial for interpreting the study" + `
    $rightQuote + "s findings accurately and recognizing the potential for " + `
    "variance in real-world applications. Further research and refined " + `
    "methodologies could help address these limitations in future analyses."
$newParagraphText = "Further research, data collection, and methodological " + `
    "refinements could help address these limitations in future analyses."

$rewordRange = $d.Content
$reworded = $rewordRange.Find.Execute($oldParagraphText, $true, $false, $false, `
    $false, $false, $true, 1, $false, $newParagraphText, 2)

if (-not $reworded) {
    Write-Host "WARNING: could not find the limitations summary paragraph to reword."
}

# --- Step 2: delete the nine bulleted limitation list items -----------------

$firstBulletText = "Only includes parcels where residential use is allowed"
$lastBulletText = "Assumes all data to be concurrent"

$startRange = $d.Content
$foundFirst = $startRange.Find.Execute($firstBulletText, $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

$endRange = $d.Content
$foundLast = $endRange.Find.Execute($lastBulletText, $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

if ($foundFirst -and $foundLast) {
    # Grow the end boundary to cover the whole last paragraph (including its
    # paragraph mark) so the entire run of list-item paragraphs disappears
    # cleanly, leaving the bookmarks that sit right after them untouched.
    $lastBulletParagraph = $endRange.Paragraphs(1)
    $bulletsRange = $d.Range($startRange.Start, $lastBulletParagraph.Range.End)
    $bulletsRange.Delete()
} else {
    Write-Host "WARNING: could not find the bulleted limitations list to delete."
}
